$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nested")

$ws.Range("C1").Value = "list#key?toMap=key"
$ws.Range("D1").Value = "list#value?toMap=value"
